{"js": "// Replace each \"NN\u00d7NN=\" arithmetic-problem cell text with its new value.\n// Source and target strings are all unique in this document, so a direct\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"24\u00d718=\", \"49\u00d738=\"],\n  [\"26\u00d765=\", \"59\u00d712=\"],\n  [\"36\u00d713=\", \"48\u00d747=\"],\n  [\"51\u00d756=\", \"62\u00d778=\"],\n  [\"59\u00d767=\", \"62\u00d749=\"],\n  [\"72\u00d758=\", \"68\u00d733=\"],\n  [\"64\u00d750=\", \"87\u00d747=\"],\n  [\"14\u00d738=\", \"67\u00d794=\"],\n  [\"87\u00d755=\", \"62\u00d751=\"],\n  [\"60\u00d799=\", \"84\u00d711=\"],\n  [\"65\u00d762=\", \"17\u00d721=\"],\n  [\"43\u00d747=\", \"73\u00d779=\"],\n  [\"47\u00d771=\", \"42\u00d718=\"],\n  [\"19\u00d784=\", \"27\u00d762=\"],\n  [\"73\u00d728=\", \"90\u00d749=\"],\n  [\"66\u00d711=\", \"47\u00d772=\"],\n  [\"78\u00d784=\", \"57\u00d779=\"],\n  [\"24\u00d727=\", \"61\u00d716=\"],\n  [\"81\u00d723=\", \"87\u00d724=\"],\n  [\"96\u00d768=\", \"24\u00d715=\"],\n  [\"19\u00d783=\", \"11\u00d764=\"],\n  [\"25\u00d755=\", \"11\u00d727=\"],\n  [\"36\u00d770=\", \"37\u00d770=\"],\n  [\"76\u00d792=\", \"80\u00d716=\"],\n  [\"14\u00d797=\", \"77\u00d716=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"NN\u00d7NN=\" arithmetic-problem cell text with its new value.\n# Source and target strings are all unique in this document, so a direct\n# Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ old = \"24\u00d718=\"; new = \"49\u00d738=\" },\n    @{ old = \"26\u00d765=\"; new = \"59\u00d712=\" },\n    @{ old = \"36\u00d713=\"; new = \"48\u00d747=\" },\n    @{ old = \"51\u00d756=\"; new = \"62\u00d778=\" },\n    @{ old = \"59\u00d767=\"; new = \"62\u00d749=\" },\n    @{ old = \"72\u00d758=\"; new = \"68\u00d733=\" },\n    @{ old = \"64\u00d750=\"; new = \"87\u00d747=\" },\n    @{ old = \"14\u00d738=\"; new = \"67\u00d794=\" },\n    @{ old = \"87\u00d755=\"; new = \"62\u00d751=\" },\n    @{ old = \"60\u00d799=\"; new = \"84\u00d711=\" },\n    @{ old = \"65\u00d762=\"; new = \"17\u00d721=\" },\n    @{ old = \"43\u00d747=\"; new = \"73\u00d779=\" },\n    @{ old = \"47\u00d771=\"; new = \"42\u00d718=\" },\n    @{ old = \"19\u00d784=\"; new = \"27\u00d762=\" },\n    @{ old = \"73\u00d728=\"; new = \"90\u00d749=\" },\n    @{ old = \"66\u00d711=\"; new = \"47\u00d772=\" },\n    @{ old = \"78\u00d784=\"; new = \"57\u00d779=\" },\n    @{ old = \"24\u00d727=\"; new = \"61\u00d716=\" },\n    @{ old = \"81\u00d723=\"; new = \"87\u00d724=\" },\n    @{ old = \"96\u00d768=\"; new = \"24\u00d715=\" },\n    @{ old = \"19\u00d783=\"; new = \"11\u00d764=\" },\n    @{ old = \"25\u00d755=\"; new = \"11\u00d727=\" },\n    @{ old = \"36\u00d770=\"; new = \"37\u00d770=\" },\n    @{ old = \"76\u00d792=\"; new = \"80\u00d716=\" },\n    @{ old = \"14\u00d797=\"; new = \"77\u00d716=\" }\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}\n"}
